$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "AgregarProducto"

# Header row
$ws.Range("A1").Value = "url"
$ws.Range("B1").Value = "producto"
$ws.Range("C1").Value = "cantidadProducto"

# Data row
$ws.Range("A2").Value = "https://es.aliexpress.com/"
$ws.Range("B2").Value = "ADDONEE-cárdigan con cremallera para hombre"
$ws.Range("C2").Value = 3

# Header style: bold, size 12
$headerRange = $ws.Range("A1:C1")
$headerFont = $headerRange.Font
$headerFont.Size = 12
$headerFont.Bold = $true
$ws.Rows.Item(1).RowHeight = 15.75

# Hyperlink on A2 pointing to the product page, displaying the base url
$ws.Hyperlinks.Add($ws.Range("A2"), "https://es.aliexpress.com/item/1005006041802486.html", "", "", "https://es.aliexpress.com/")

# Column widths (approximate best-fit)
$ws.Columns.Item(1).ColumnWidth = 24.02
$ws.Columns.Item(2).ColumnWidth = 43.02
$ws.Columns.Item(3).ColumnWidth = 18.02

# Page setup
$ws.PageSetup.Orientation = 1

# Selection matching target
$null = $ws.Range("F8").Select()
